# Auto-generated Excel COM-interop script
# Applies market-price / profit-column refresh values to the Jenova_Profits workbook
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW and WVR leve-profit tables.

$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2128.6128
$ws.Range("J17").Value = 2156.2334
$ws.Range("L17").Value = 6468.7002
$ws.Range("N17").Value = -6804.7002
$ws.Range("H62").Value = 6948730
$ws.Range("I62").Value = 12502066
$ws.Range("K62").Value = 12502066
$ws.Range("M62").Value = -12501442
$ws.Range("H65").Value = 6948730
$ws.Range("I65").Value = 12502066
$ws.Range("K65").Value = 62510330
$ws.Range("M65").Value = -62507210
$ws.Range("H70").Value = 253413.25
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 253413.25
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 760239.75
$ws.Range("M70").ClearContents()   # was -1228.5
$ws.Range("N70").Value = -760779.75
$ws.Range("H73").Value = 253413.25
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 253413.25
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 760239.75
$ws.Range("M73").ClearContents()   # was -562.5
$ws.Range("N73").Value = -762111.75
$ws.Range("H88").Value = 3560.8333
$ws.Range("I88").Value = 3316.5
$ws.Range("J88").Value = 3683
$ws.Range("K88").Value = 3316.5
$ws.Range("L88").Value = 3683
$ws.Range("M88").Value = -2910.5
$ws.Range("N88").Value = -4495
$ws.Range("H91").Value = 3560.8333
$ws.Range("I91").Value = 3316.5
$ws.Range("J91").Value = 3683
$ws.Range("K91").Value = 3316.5
$ws.Range("L91").Value = 3683
$ws.Range("M91").Value = -1912.5
$ws.Range("N91").Value = -6491
$ws.Range("H131").Value = 5335.7144
$ws.Range("I131").Value = 4270
$ws.Range("K131").Value = 12810
$ws.Range("M131").Value = -7770
$ws.Range("H137").Value = 1182526.2
$ws.Range("I137").Value = 835245.25
$ws.Range("K137").Value = 2505735.75
$ws.Range("M137").Value = -2503185.75
$ws.Range("H138").Value = 6738.9624
$ws.Range("J138").Value = 9669.588
$ws.Range("L138").Value = 29008.764
$ws.Range("N138").Value = -39288.764

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H21").Value = 3699.5
$ws.Range("I21").Value = 3699.5
$ws.Range("K21").Value = 3699.5
$ws.Range("M21").Value = -3325.5
$ws.Range("H32").Value = 5088.0625
$ws.Range("I32").Value = 4345.2554
$ws.Range("K32").Value = 4345.2554
$ws.Range("M32").Value = -4058.2554
$ws.Range("H97").Value = 1699.8422
$ws.Range("I97").Value = 2037.6666
$ws.Range("J97").Value = 433
$ws.Range("K97").Value = 2037.6666
$ws.Range("L97").Value = 433
$ws.Range("M97").Value = -1541.6666
$ws.Range("N97").Value = -1425
$ws.Range("H132").Value = 3902.851
$ws.Range("I132").Value = 2872.5945
$ws.Range("K132").Value = 8617.783500000001
$ws.Range("M132").Value = -6087.783500000001

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2213.7727
$ws.Range("I20").Value = 1579.5385
$ws.Range("K20").Value = 1579.5385
$ws.Range("M20").Value = -1332.5385
$ws.Range("H94").Value = 868
$ws.Range("I94").Value = 957.7857
$ws.Range("K94").Value = 957.7857
$ws.Range("M94").Value = -506.7857
$ws.Range("H99").Value = 4103.5713
$ws.Range("I99").Value = 3492.8
$ws.Range("K99").Value = 3492.8
$ws.Range("M99").Value = -1994.8

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 36352.25
$ws.Range("I31").Value = 1962
$ws.Range("K31").Value = 1962
$ws.Range("M31").Value = -1667
$ws.Range("H34").Value = 36352.25
$ws.Range("I34").Value = 1962
$ws.Range("K34").Value = 1962
$ws.Range("M34").Value = -1760
$ws.Range("H134").Value = 1337125.6
$ws.Range("I134").Value = 1004334.7
$ws.Range("K134").Value = 3013004.1
$ws.Range("M134").Value = -3010469.1

# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 812.5
$ws.Range("I64").Value = 812.5
$ws.Range("K64").Value = 2437.5
$ws.Range("M64").Value = -2167.5
$ws.Range("H67").Value = 812.5
$ws.Range("I67").Value = 812.5
$ws.Range("K67").Value = 2437.5
$ws.Range("M67").Value = -1501.5
$ws.Range("H82").Value = 7791.6665
$ws.Range("J82").Value = 7350
$ws.Range("L82").Value = 22050
$ws.Range("N82").Value = -22862
$ws.Range("H85").Value = 7791.6665
$ws.Range("J85").Value = 7350
$ws.Range("L85").Value = 22050
$ws.Range("N85").Value = -24858
$ws.Range("H136").Value = 5060.5713
$ws.Range("I136").Value = 3085.2
$ws.Range("K136").Value = 9255.599999999999
$ws.Range("M136").Value = -4155.599999999999
$ws.Range("H137").Value = 2033.409
$ws.Range("I137").Value = 1196.579
$ws.Range("J137").Value = 7333.3335
$ws.Range("K137").Value = 3589.737
$ws.Range("L137").Value = 22000.0005
$ws.Range("M137").Value = 1510.263
$ws.Range("N137").Value = -32200.0005

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 974.4516
$ws.Range("I97").Value = 902.7895
$ws.Range("J97").Value = 1087.9166
$ws.Range("K97").Value = 902.7895
$ws.Range("L97").Value = 1087.9166
$ws.Range("M97").Value = -406.7895
$ws.Range("N97").Value = -2079.9166
$ws.Range("H102").Value = 2302.1082
$ws.Range("I102").Value = 1497.6786
$ws.Range("K102").Value = 1497.6786
$ws.Range("M102").Value = 124.3214
$ws.Range("H131").Value = 34163
$ws.Range("J131").Value = 34163
$ws.Range("L131").Value = 34163
$ws.Range("N131").Value = -44243
$ws.Range("H132").Value = 357957.47
$ws.Range("I132").Value = 479665.8
$ws.Range("K132").Value = 1438997.4
$ws.Range("M132").Value = -1436467.4

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 1252125.4
$ws.Range("I40").Value = 1614631.2
$ws.Range("K40").Value = 1614631.2
$ws.Range("M40").Value = -1614495.2
$ws.Range("H93").Value = 3127.4443
$ws.Range("I93").Value = 2400
$ws.Range("J93").Value = 3491.1667
$ws.Range("K93").Value = 2400
$ws.Range("L93").Value = 3491.1667
$ws.Range("M93").Value = -1152
$ws.Range("N93").Value = -5987.1667
$ws.Range("H132").Value = 4617.778
$ws.Range("I132").Value = 2640
$ws.Range("K132").Value = 7920
$ws.Range("M132").Value = -5390

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1993.7646
$ws.Range("I81").Value = 1524.3
$ws.Range("J81").Value = 2664.4285
$ws.Range("K81").Value = 3048.6
$ws.Range("L81").Value = 5328.857
$ws.Range("M81").Value = -1987.6
$ws.Range("N81").Value = -7450.857
$ws.Range("H84").Value = 1993.7646
$ws.Range("I84").Value = 1524.3
$ws.Range("J84").Value = 2664.4285
$ws.Range("K84").Value = 15243
$ws.Range("L84").Value = 26644.285
$ws.Range("M84").Value = -9939
$ws.Range("N84").Value = -37252.285
$ws.Range("H122").Value = 40004852
$ws.Range("I122").Value = 50004370
$ws.Range("K122").Value = 150013110
$ws.Range("M122").Value = -150010660
$ws.Range("H132").Value = 68851.25
$ws.Range("I132").Value = 7174.25
$ws.Range("K132").Value = 21522.75
$ws.Range("M132").Value = -18992.75
